$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 106, shifting existing rows 106:147 down to 107:148.
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row 106 with the new data point.
$ws.Cells.Item(106, 1).Value = 5
$ws.Cells.Item(106, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(106, 3).Value = "Maule"
$ws.Cells.Item(106, 4).Value = 44992
$ws.Cells.Item(106, 5).Value = 7
$ws.Cells.Item(106, 6).Value = 100112001
$ws.Cells.Item(106, 7).Value = "Berenjena"
$ws.Cells.Item(106, 8).Value = "Sin especificar"
$ws.Cells.Item(106, 9).Value = "Primera"
$ws.Cells.Item(106, 10).Value = 150
$ws.Cells.Item(106, 11).Value = 8000
$ws.Cells.Item(106, 12).Value = 8000
$ws.Cells.Item(106, 13).Value = 8000
$ws.Cells.Item(106, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(106, 15).Value = "Región del Maule"
$ws.Cells.Item(106, 16).Value = 160
$ws.Cells.Item(106, 17).Value = 50
$ws.Cells.Item(106, 18).Value = "Hortaliza"
